$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "49.992.56"
$ws.Range("E2").Value = "  +4.02%  "

# Row 3
$ws.Range("D3").Value = "2.647.79"
$ws.Range("E3").Value = "  +6.05%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "113.92"
$ws.Range("E5").Value = "  +8.07%  "

# Row 6
$ws.Range("D6").Value = "326.95"
$ws.Range("E6").Value = "  +2.74%  "

# Row 7
$ws.Range("E7").Value = "  +1.75%  "

# Row 8
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").Value = "0.557"
$ws.Range("E9").Value = "  +3.99%  "

# Row 10
$ws.Range("D10").Value = "41.00"
$ws.Range("E10").Value = "  +5.81%  "

# Row 11
$ws.Range("D11").Value = "20.25"
$ws.Range("E11").Value = "  +0.24%  "

# Row 12
$ws.Range("E12").Value = "  +2.63%  "

# Row 13
$ws.Range("E13").Value = "  +0.98%  "

# Row 14
$ws.Range("D14").Value = "7.39"
$ws.Range("E14").Value = "  +4.60%  "

# Row 15
$ws.Range("D15").Value = "3.063.60"
$ws.Range("E15").Value = "  +6.12%  "

# Row 16
$ws.Range("D16").Value = "2.647.94"
$ws.Range("E16").Value = "  +5.83%  "

# Row 17
$ws.Range("E17").Value = "  +5.57%  "

# Row 18
$ws.Range("D18").Value = "49.914.31"
$ws.Range("E18").Value = "  +4.21%  "

# Row 19
$ws.Range("D19").Value = "13.22"
$ws.Range("E19").Value = "  +2.72%  "

# Row 20
$ws.Range("E20").Value = "  +2.71%  "

# Row 21
$ws.Range("E21").Value = "  -1.85%  "

# Row 22
$ws.Range("E22").Value = "  +3.43%  "

# Row 23
$ws.Range("D23").Value = "72.09"
$ws.Range("E23").Value = "  +1.51%  "

# Row 24
$ws.Range("D24").Value = "276.77"
$ws.Range("E24").Value = "  +2.55%  "

# Row 25
$ws.Range("D25").Value = "2.60"
$ws.Range("E25").Value = "  +3.25%  "

# Row 26
$ws.Range("E26").Value = "  +4.01%  "

# Row 27
$ws.Range("E27").Value = "  -0.05%  "

# Row 28
$ws.Range("D28").Value = "10.00"
$ws.Range("E28").Value = "  +3.16%  "

# Row 29
$ws.Range("E29").Value = "  -2.27%  "

# Row 30
$ws.Range("D30").Value = "36.06"
$ws.Range("E30").Value = "  +4.79%  "

# Row 31
$ws.Range("E31").Value = "  +2.21%  "

# Row 32
$ws.Range("D32").Value = "50.34"
$ws.Range("E32").Value = "  +2.10%  "

# Row 33
$ws.Range("E33").Value = "  +3.33%  "

# Row 34
$ws.Range("D34").Value = "19.47"
$ws.Range("E34").Value = "  +3.10%  "

# Row 35
$ws.Range("E35").Value = "  +5.11%  "

# Row 36
$ws.Range("E36").Value = "  -0.04%  "

# Row 37
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "4.97"
$ws.Range("E37").Value = "  +8.68%  "

# Row 38
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "2.07"
$ws.Range("E38").Value = "  +7.00%  "

# Row 39
$ws.Range("D39").Value = "3.10"
$ws.Range("E39").Value = "  +8.30%  "

# Row 40
$ws.Range("D40").Value = "123.85"
$ws.Range("E40").Value = "  +1.53%  "

# Row 41
$ws.Range("E41").Value = "  +2.20%  "

# Row 42
$ws.Range("E42").Value = "  +0.48%  "

# Row 43
$ws.Range("D43").Value = "22.10"
$ws.Range("E43").Value = "  -1.61%  "

# Row 44
$ws.Range("E44").Value = "  +4.72%  "

# Row 45
$ws.Range("D45").Value = "2.085.87"
$ws.Range("E45").Value = "  +4.43%  "

# Row 46
$ws.Range("E46").Value = "  +6.88%  "

# Row 47
$ws.Range("E47").Value = "  +16.57%  "

# Row 48
$ws.Range("E48").Value = "  +6.42%  "

# Row 49
$ws.Range("E49").Value = "  +2.79%  "

# Row 50
$ws.Range("D50").Value = "5.41"
$ws.Range("E50").Value = "  +5.06%  "

# Row 51
$ws.Range("D51").Value = "59.89"
$ws.Range("E51").Value = "  +6.35%  "
